# Lesson 8.2 Deployment Infrastructure — title slide edit
# Commit: "Remove CS 5500 from titles"
#
# The title placeholder on slide 1 originally had two paragraphs:
#   "CS 4530 & CS 5500"
#   "Software Engineering"
# It becomes a single paragraph:
#   "CS 4530 Software Engineering"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the title shape (placeholder "CS 4530 & CS 5500...") on slide 1
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $text = $shape.TextFrame.TextRange.Text
        if ($text -like "CS 4530*CS 5500*Software Engineering*") {
            $shape.TextFrame.TextRange.Text = "CS 4530 Software Engineering"
        }
    }
}
